$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.631.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.061.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.670'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.82'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '58.81'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.363'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("E12").Value = '  -3.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.935'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.362.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.026.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.553.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0862'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.52%  '
$ws.Range("E26").Value = '  +1.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.123'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("E31").Value = '  +8.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.08%  '
$ws.Range("E33").Value = '  -4.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0601'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.82%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.22'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0822'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("E39").Value = '  -3.95%  '
$ws.Range("E40").Value = '  -5.65%  '
$ws.Range("E41").Value = '  -1.96%  '
$ws.Range("E42").Value = '  -7.76%  '
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0911'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.417.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.82%  '
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.247.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.16%  '
